$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 787.6
$ws.Range("I86").Value = 489
$ws.Range("J86").Value = 986.6667
$ws.Range("K86").Value = 489
$ws.Range("L86").Value = 986.6667
$ws.Range("M86").Value = 634
$ws.Range("N86").Value = -3232.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 787.6
$ws.Range("I89").Value = 489
$ws.Range("J89").Value = 986.6667
$ws.Range("K89").Value = 2445
$ws.Range("L89").Value = 4933.3335
$ws.Range("M89").Value = 3171
$ws.Range("N89").Value = -16165.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2369.6
$ws.Range("I132").Value = 2974
$ws.Range("K132").Value = 8922
$ws.Range("M132").Value = -6392

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1503.2
$ws.Range("I137").Value = 1379
$ws.Range("K137").Value = 4137
$ws.Range("M137").Value = -1587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1429.1
$ws.Range("I2").Value = 562.25
$ws.Range("K2").Value = 562.25
$ws.Range("M2").Value = -449.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 750
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1429.1
$ws.Range("I116").Value = 562.25
$ws.Range("K116").Value = 562.25
$ws.Range("M116").Value = 1731.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1275.8823
$ws.Range("I132").Value = 1323.8462
$ws.Range("K132").Value = 3971.5386
$ws.Range("M132").Value = -1441.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 750
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1429.1
$ws.Range("I3").Value = 562.25
$ws.Range("K3").Value = 562.25
$ws.Range("M3").Value = -448.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 29474.75
$ws.Range("I97").Value = 17966.666
$ws.Range("K97").Value = 17966.666
$ws.Range("M97").Value = -16975.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2479.4644
$ws.Range("I134").Value = 1469.7142
$ws.Range("J134").Value = 3489.2144
$ws.Range("K134").Value = 4409.142599999999
$ws.Range("L134").Value = 10467.6432
$ws.Range("M134").Value = -1874.142599999999
$ws.Range("N134").Value = -15537.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6000
$ws.Range("J4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("N4").Value = -8224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 102.46154
$ws.Range("I7").Value = 49.333332
$ws.Range("K7").Value = 49.333332
$ws.Range("M7").Value = 63.666668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3363.4583
$ws.Range("I31").Value = 2349.0588
$ws.Range("J31").Value = 5827
$ws.Range("K31").Value = 2349.0588
$ws.Range("L31").Value = 5827
$ws.Range("M31").Value = -2054.0588
$ws.Range("N31").Value = -6417

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3363.4583
$ws.Range("I34").Value = 2349.0588
$ws.Range("J34").Value = 5827
$ws.Range("K34").Value = 2349.0588
$ws.Range("L34").Value = 5827
$ws.Range("M34").Value = -2147.0588
$ws.Range("N34").Value = -6231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1889.0667
$ws.Range("I58").Value = 1121
$ws.Range("K58").Value = 1121
$ws.Range("M58").Value = -918

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9800
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 9800
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 32011.445
$ws.Range("J95").Value = 32011.445
$ws.Range("L95").Value = 32011.445
$ws.Range("N95").Value = -37503.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 789.8182
$ws.Range("I107").Value = 786.375
$ws.Range("K107").Value = 786.375
$ws.Range("M107").Value = 1133.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 983
$ws.Range("I132").Value = 983
$ws.Range("K132").Value = 2949
$ws.Range("M132").Value = -419

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2415.158
$ws.Range("J134").Value = 3745
$ws.Range("L134").Value = 11235
$ws.Range("N134").Value = -16305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1889.0667
$ws.Range("I136").Value = 1121
$ws.Range("K136").Value = 3363
$ws.Range("M136").Value = -813

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10488655
$ws.Range("J4").Value = 802
$ws.Range("L4").Value = 2406
$ws.Range("N4").Value = -2630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 164.33333
$ws.Range("I103").Value = 137
$ws.Range("K103").Value = 411
$ws.Range("M103").Value = 468

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -4060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1499.5
$ws.Range("I9").Value = 1499.5
$ws.Range("K9").Value = 1499.5
$ws.Range("M9").Value = -1329.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8999.25
$ws.Range("J80").Value = 8999.25
$ws.Range("L80").Value = 8999.25
$ws.Range("N80").Value = -10995.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8999.25
$ws.Range("J83").Value = 8999.25
$ws.Range("L83").Value = 44996.25
$ws.Range("N83").Value = -54980.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5631.636
$ws.Range("I22").Value = 3494.4443
$ws.Range("J22").Value = 15249
$ws.Range("K22").Value = 3494.4443
$ws.Range("L22").Value = 15249
$ws.Range("M22").Value = -3199.4443
$ws.Range("N22").Value = -15839

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5631.636
$ws.Range("I27").Value = 3494.4443
$ws.Range("J27").Value = 15249
$ws.Range("K27").Value = 3494.4443
$ws.Range("L27").Value = 15249
$ws.Range("M27").Value = -3387.4443
$ws.Range("N27").Value = -15463

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 943.0769
$ws.Range("I55").Value = 661.75
$ws.Range("J55").Value = 1393.2
$ws.Range("K55").Value = 661.75
$ws.Range("L55").Value = 1393.2
$ws.Range("M55").Value = -488.75
$ws.Range("N55").Value = -1739.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4399.091
$ws.Range("I61").Value = 4339
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4339
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -4137
$ws.Range("N61").Value = -5404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2914.7144
$ws.Range("I68").Value = 2100
$ws.Range("J68").Value = 4951.5
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 4951.5
$ws.Range("M68").Value = -1351
$ws.Range("N68").Value = -6449.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2914.7144
$ws.Range("I71").Value = 2100
$ws.Range("J71").Value = 4951.5
$ws.Range("K71").Value = 10500
$ws.Range("L71").Value = 24757.5
$ws.Range("M71").Value = -6756
$ws.Range("N71").Value = -32245.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2815.2222
$ws.Range("I82").Value = 3479
$ws.Range("J82").Value = 2483.3333
$ws.Range("K82").Value = 3479
$ws.Range("L82").Value = 2483.3333
$ws.Range("M82").Value = -3118
$ws.Range("N82").Value = -3205.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2815.2222
$ws.Range("I85").Value = 3479
$ws.Range("J85").Value = 2483.3333
$ws.Range("K85").Value = 3479
$ws.Range("L85").Value = 2483.3333
$ws.Range("M85").Value = -2231
$ws.Range("N85").Value = -4979.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4399.091
$ws.Range("I113").Value = 4339
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4339
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2169
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2185.1714
$ws.Range("I132").Value = 1951.129
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 5853.387
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3323.387
$ws.Range("N132").Value = -17057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1249.5
$ws.Range("I17").Value = 1249.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 1249.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1077.5
$ws.Range("N17").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2955.5715
$ws.Range("I126").Value = 2731.5
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 8194.5
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -5724.5
$ws.Range("N126").Value = -17840
